$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing column C. This shifts the
# original column C (and its column-width formatting) two positions to
# the right, landing on column E.
$ws.Range("C:D").Insert()

# Row 1 holds the "as of" date labels for each data column. The newest
# date used to live in B1 ("Jun_13"); now a newer date ("Jun_17") takes
# its place and the old value moves one column over into the newly
# inserted column D. A brand new label ("Jun_15") is put in C1.
$ws.Cells.Item(1, 2).Value = "Jun_17"
$ws.Cells.Item(1, 3).Value = "Jun_15"
$ws.Cells.Item(1, 4).Value = "Jun_13"
# E1 already holds the original column C value ("Jun_10") after the insert.

# The two newly inserted data columns (C and D) get the same default
# "UN" placeholder value that column B already uses, for every data row.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Cosmetic: give all three date columns the same fixed width (8
# characters) that the original data column already used. (The insert
# above shifted the original column C -- and its custom-width flag --
# to column E, so re-assert it explicitly there too.)
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(4).ColumnWidth = 7.166666666666667
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667
